$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Scenarios")

$ws.Range("E6").Value = "Global, MissingParam"
$ws.Range("A6").Value = "TestScenario_missingParam"
$ws.Range("B6").Value = "Indiv1"
$ws.Range("F6").Value = "Aciclovir_iv_250mg"
$ws.Range("G6").Value = "0, 24, 60"
$ws.Range("H6").Value = "h"
$ws.Range("L6").Value = "Aciclovir.pkml"
